# Insert a new weekly price-report row for Espinaca (Femacal de La Calera)
# right above what is currently row 622. All rows from the old 622 down to
# 653 shift down by one (to 623..654), and the new row 622 receives the
# latest observation. Final used range becomes A1:R654.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 622 (and everything below it) down by one row.
$ws.Rows.Item(622).Insert()

# Populate the newly inserted row 622 with the new weekly record.
$ws.Cells.Item(622, 1).Value = 3
$ws.Cells.Item(622, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(622, 3).Value = "Coquimbo"
$ws.Cells.Item(622, 4).Value = 45267
$ws.Cells.Item(622, 5).Value = 5
$ws.Cells.Item(622, 6).Value = 100112012
$ws.Cells.Item(622, 7).Value = "Espinaca"
$ws.Cells.Item(622, 8).Value = "Sin especificar"
$ws.Cells.Item(622, 9).Value = "Primera"
$ws.Cells.Item(622, 10).Value = 100
$ws.Cells.Item(622, 11).Value = 6000
$ws.Cells.Item(622, 12).Value = 7000
$ws.Cells.Item(622, 13).Value = 6500
$ws.Cells.Item(622, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(622, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(622, 16).Value = 2167
$ws.Cells.Item(622, 17).Value = 3
$ws.Cells.Item(622, 18).Value = "Hortaliza"
